# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handback packages have come back "in sync with en-US":
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears
#     (Overview sheet zh-cn/de-de columns, and the Status column on the
#     per-language sheets).
#   - The "Latest Target File" (I) and "Latest Handback File" (J) columns on
#     the zh-cn/de-de sheets are now populated with the handed-back xliff
#     file info, and I2/I3 become "a.md" hyperlinks (matching A2's link).
#   - The "Latest Handback DateTime" (K) column is stamped with the handback
#     time for each language.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6901a4226f0a2e0e4e37f5fab096fe9758bcdd7c/e2e/a.md"
$bMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6901a4226f0a2e0e4e37f5fab096fe9758bcdd7c/e2e/b.md"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn (E) / de-de (F) status cells for both rows
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Columns E and F widen to fit the longer status text.
$overview.Range("E1").ColumnWidth = 29.1
$overview.Range("F1").ColumnWidth = 29.1

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("I2").Value = "a.md"
$zhcn.Range("I2").Style = "Hyperlink"
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-26 20:36:21"

$zhcn.Range("I3").Value = "a.md"
$zhcn.Range("I3").Style = "Hyperlink"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-26 20:36:21"

# Rebuild the hyperlinks collection so the new I2/I3 "a.md" links sit
# alongside the existing A2/A3 links, in the same relative order as the
# final sheet (A2, I2, A3, I3).
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $bMdUrl, [Type]::Missing, [Type]::Missing, "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")

# Columns C (Status) and J (Latest Handback File) widen for the new text.
$zhcn.Range("C1").ColumnWidth = 29.1
$zhcn.Range("J1").ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = "a.md"
$dede.Range("I2").Style = "Hyperlink"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-26 20:36:28"

$dede.Range("I3").Value = "a.md"
$dede.Range("I3").Style = "Hyperlink"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-26 20:36:28"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$dede.Hyperlinks.Add($dede.Range("I2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")
$dede.Hyperlinks.Add($dede.Range("A3"), $bMdUrl, [Type]::Missing, [Type]::Missing, "b.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md")

# Columns C (Status) and J (Latest Handback File) widen for the new text.
$dede.Range("C1").ColumnWidth = 29.1
$dede.Range("J1").ColumnWidth = 39.1
